# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (period) rows for worker 84062118 (LUIS JAVIER GAZABON PEREZ)
# are resequenced so period 2208 becomes the first row, and a new row is
# inserted right after it for worker 22801531 (YESENIA DEL MAR RODRIGUEZ
# IRIARTE), also for period 2208 but with a different "Valor Mora" (1160000).
# All the remaining periods for LUIS JAVIER GAZABON PEREZ shift up by one row
# and are renumbered in ascending order, ending with period 2307 (which keeps
# its special "Salario Basico" value of 1333) on the last data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: CC / 84062118 / LUIS JAVIER GAZABON PEREZ / 2208 / 40000 / 1000000
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "84062118"
$ws.Range("D16").Value = "LUIS JAVIER GAZABON PEREZ"
$ws.Range("E16").Value = "2208"
$ws.Range("F16").Value = 40000
$ws.Range("G16").Value = 1000000

# Row 17: CC / 22801531 / YESENIA DEL MAR RODRIGUEZ IRIARTE / 2208 / 40000 / 1160000
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "22801531"
$ws.Range("D17").Value = "YESENIA DEL MAR RODRIGUEZ IRIARTE"
$ws.Range("E17").Value = "2208"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1160000

# Row 18: CC / 84062118 / LUIS JAVIER GAZABON PEREZ / 2209 / 40000 / 1000000
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "84062118"
$ws.Range("D18").Value = "LUIS JAVIER GAZABON PEREZ"
$ws.Range("E18").Value = "2209"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1000000

# Row 19: CC / 84062118 / LUIS JAVIER GAZABON PEREZ / 2210 / 40000 / 1000000
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "84062118"
$ws.Range("D19").Value = "LUIS JAVIER GAZABON PEREZ"
$ws.Range("E19").Value = "2210"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 1000000

# Row 20: CC / 84062118 / LUIS JAVIER GAZABON PEREZ / 2211 / 40000 / 1000000
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "84062118"
$ws.Range("D20").Value = "LUIS JAVIER GAZABON PEREZ"
$ws.Range("E20").Value = "2211"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 1000000

# Row 21: CC / 84062118 / LUIS JAVIER GAZABON PEREZ / 2212 / 40000 / 1000000
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "84062118"
$ws.Range("D21").Value = "LUIS JAVIER GAZABON PEREZ"
$ws.Range("E21").Value = "2212"
$ws.Range("F21").Value = 40000
$ws.Range("G21").Value = 1000000

# Row 22: CC / 84062118 / LUIS JAVIER GAZABON PEREZ / 2301 / 40000 / 1000000
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "84062118"
$ws.Range("D22").Value = "LUIS JAVIER GAZABON PEREZ"
$ws.Range("E22").Value = "2301"
$ws.Range("F22").Value = 40000
$ws.Range("G22").Value = 1000000

# Row 23: CC / 84062118 / LUIS JAVIER GAZABON PEREZ / 2302 / 40000 / 1000000
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "84062118"
$ws.Range("D23").Value = "LUIS JAVIER GAZABON PEREZ"
$ws.Range("E23").Value = "2302"
$ws.Range("F23").Value = 40000
$ws.Range("G23").Value = 1000000

# Row 24: CC / 84062118 / LUIS JAVIER GAZABON PEREZ / 2303 / 40000 / 1000000
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "84062118"
$ws.Range("D24").Value = "LUIS JAVIER GAZABON PEREZ"
$ws.Range("E24").Value = "2303"
$ws.Range("F24").Value = 40000
$ws.Range("G24").Value = 1000000

# Row 25: CC / 84062118 / LUIS JAVIER GAZABON PEREZ / 2304 / 40000 / 1000000
$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "84062118"
$ws.Range("D25").Value = "LUIS JAVIER GAZABON PEREZ"
$ws.Range("E25").Value = "2304"
$ws.Range("F25").Value = 40000
$ws.Range("G25").Value = 1000000

# Row 26: CC / 84062118 / LUIS JAVIER GAZABON PEREZ / 2305 / 40000 / 1000000
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "84062118"
$ws.Range("D26").Value = "LUIS JAVIER GAZABON PEREZ"
$ws.Range("E26").Value = "2305"
$ws.Range("F26").Value = 40000
$ws.Range("G26").Value = 1000000

# Row 27: CC / 84062118 / LUIS JAVIER GAZABON PEREZ / 2306 / 40000 / 1000000
$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "84062118"
$ws.Range("D27").Value = "LUIS JAVIER GAZABON PEREZ"
$ws.Range("E27").Value = "2306"
$ws.Range("F27").Value = 40000
$ws.Range("G27").Value = 1000000

# Row 28: CC / 84062118 / LUIS JAVIER GAZABON PEREZ / 2307 / 1333 / 1000000
$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "84062118"
$ws.Range("D28").Value = "LUIS JAVIER GAZABON PEREZ"
$ws.Range("E28").Value = "2307"
$ws.Range("F28").Value = 1333
$ws.Range("G28").Value = 1000000
